# Automatico via Actualizar 02-04-2021 19-36-17
# This script appends a fresh block of 14 availability-check rows
# (rows 100-113) to Sheet1, mirroring the existing repeating
# Nombre/URL/Disponibilidad/Fecha table, and refreshes the timestamp
# of the previous block (rows 86-99) to the latest check.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the timestamp on the previous run's last block (rows 86-99) ---
$prevTimestamp = 44231.7954509375
foreach ($r in 86..99) {
    $ws.Range("D$r").Value = $prevTimestamp
}

# --- Append the new block for this run (rows 100-113) ---
$newTimestamp = 44231.8165548469

# Row 100
$ws.Range("A100").Value = 'Odoo'
$ws.Range("B100").Value = 'https://www.dataintelligence-group.com/'
$ws.Range("C100").Value = 'Disponible'
$ws.Range("D100").Value = $newTimestamp
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B100"), 'https://www.dataintelligence-group.com/') | Out-Null
$ws.Range("B100").Style = "Hyperlink"

# Row 101
$ws.Range("A101").Value = 'Blackbox'
$ws.Range("B101").Value = 'https://serviciodashboard.azurewebsites.net/'
$ws.Range("C101").Value = 'Disponible'
$ws.Range("D101").Value = $newTimestamp
$ws.Range("D101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B101"), 'https://serviciodashboard.azurewebsites.net/') | Out-Null
$ws.Range("B101").Style = "Hyperlink"

# Row 102
$ws.Range("A102").Value = 'PowerBI'
$ws.Range("B102").Value = 'https://powerbi.microsoft.com/es-es/'
$ws.Range("C102").Value = 'Disponible'
$ws.Range("D102").Value = $newTimestamp
$ws.Range("D102").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B102"), 'https://powerbi.microsoft.com/es-es/') | Out-Null
$ws.Range("B102").Style = "Hyperlink"

# Row 103
$ws.Range("A103").Value = 'Dropbox'
$ws.Range("B103").Value = 'https://www.dropbox.com/'
$ws.Range("C103").Value = 'Disponible'
$ws.Range("D103").Value = $newTimestamp
$ws.Range("D103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B103"), 'https://www.dropbox.com/') | Out-Null
$ws.Range("B103").Style = "Hyperlink"

# Row 104
$ws.Range("A104").Value = 'Odoo'
$ws.Range("B104").Value = 'https://dataintelligence.store/'
$ws.Range("C104").Value = 'Disponible'
$ws.Range("D104").Value = $newTimestamp
$ws.Range("D104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B104"), 'https://dataintelligence.store/') | Out-Null
$ws.Range("B104").Style = "Hyperlink"

# Row 105
$ws.Range("A105").Value = 'GEE'
$ws.Range("B105").Value = 'https://app-data-i.users.earthengine.app/'
$ws.Range("C105").Value = 'Disponible'
$ws.Range("D105").Value = $newTimestamp
$ws.Range("D105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B105"), 'https://app-data-i.users.earthengine.app/') | Out-Null
$ws.Range("B105").Style = "Hyperlink"

# Row 106
$ws.Range("A106").Value = 'UtilidadesOdoo'
$ws.Range("B106").Value = 'https://odooutil.azurewebsites.net/'
$ws.Range("C106").Value = 'Disponible'
$ws.Range("D106").Value = $newTimestamp
$ws.Range("D106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B106"), 'https://odooutil.azurewebsites.net/') | Out-Null
$ws.Range("B106").Style = "Hyperlink"

# Row 107
$ws.Range("A107").Value = 'Filtros Dashboard'
$ws.Range("B107").Value = 'https://filtradordashboard.azurewebsites.net/'
$ws.Range("C107").Value = 'Disponible'
$ws.Range("D107").Value = $newTimestamp
$ws.Range("D107").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B107"), 'https://filtradordashboard.azurewebsites.net/') | Out-Null
$ws.Range("B107").Style = "Hyperlink"

# Row 108
$ws.Range("A108").Value = 'MapStore'
$ws.Range("B108").Value = 'https://ide.dataintelligence-group.com/mapstore/#/'
$ws.Range("C108").Value = 'Disponible'
$ws.Range("D108").Value = $newTimestamp
$ws.Range("D108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B108"), 'https://ide.dataintelligence-group.com/mapstore/', "/") | Out-Null
$ws.Range("B108").Style = "Hyperlink"

# Row 109
$ws.Range("A109").Value = 'GeoServer'
$ws.Range("B109").Value = 'https://ide.dataintelligence-group.com/geoserver/web/?0'
$ws.Range("C109").Value = 'Disponible'
$ws.Range("D109").Value = $newTimestamp
$ws.Range("D109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B109"), 'https://ide.dataintelligence-group.com/geoserver/web/?0') | Out-Null
$ws.Range("B109").Style = "Hyperlink"

# Row 110
$ws.Range("A110").Value = 'Tomcat'
$ws.Range("B110").Value = 'https://ide.dataintelligence-group.com/'
$ws.Range("C110").Value = 'Disponible'
$ws.Range("D110").Value = $newTimestamp
$ws.Range("D110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B110"), 'https://ide.dataintelligence-group.com/') | Out-Null
$ws.Range("B110").Style = "Hyperlink"

# Row 111
$ws.Range("A111").Value = 'Shiny'
$ws.Range("B111").Value = 'https://rpubs.com/dataintelligence/'
$ws.Range("C111").Value = 'Disponible'
$ws.Range("D111").Value = $newTimestamp
$ws.Range("D111").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B111"), 'https://rpubs.com/dataintelligence/') | Out-Null
$ws.Range("B111").Style = "Hyperlink"

# Row 112
$ws.Range("A112").Value = 'Github'
$ws.Range("B112").Value = 'https://github.com/Sud-Austral/'
$ws.Range("C112").Value = 'Disponible'
$ws.Range("D112").Value = $newTimestamp
$ws.Range("D112").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B112"), 'https://github.com/Sud-Austral/') | Out-Null
$ws.Range("B112").Style = "Hyperlink"

# Row 113
$ws.Range("A113").Value = 'EZ Exporter'
$ws.Range("B113").Value = 'https://ezexporter.highviewapps.com/exports/export-profile/'
$ws.Range("C113").Value = 'Disponible'
$ws.Range("D113").Value = $newTimestamp
$ws.Range("D113").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Hyperlinks.Add($ws.Range("B113"), 'https://ezexporter.highviewapps.com/exports/export-profile/') | Out-Null
$ws.Range("B113").Style = "Hyperlink"

Write-Output "Added rows 100-113 and refreshed timestamps for rows 86-99."

